$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on price cells whose new values would otherwise
# be auto-coerced to numeric by Excel (the source data stores these as text).
$textCells = @('D5', 'D6', 'D7', 'D8', 'D10', 'D11', 'D12', 'D13', 'D14', 'D16', 'D18', 'D20', 'D22', 'D23', 'D24', 'D25', 'D26', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D39', 'D41', 'D42', 'D43', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Update coin price / 1h-volume figures (and the two reordered rows).
$ws.Range("D2").Value = '42.200.69'
$ws.Range("E2").Value = '  +0.79%  '
$ws.Range("D3").Value = '2.286.10'
$ws.Range("E3").Value = '  +0.40%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '154.47'
$ws.Range("E5").Value = '  +15,332.58%  '
$ws.Range("D6").Value = '305.96'
$ws.Range("E6").Value = '  +0.70%  '
$ws.Range("D7").Value = '94.81'
$ws.Range("E7").Value = '  +1.74%  '
$ws.Range("D8").Value = '0.532'
$ws.Range("E8").Value = '  +0.17%  '
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").Value = '0.493'
$ws.Range("E10").Value = '  +1.10%  '
$ws.Range("D11").Value = '34.39'
$ws.Range("E11").Value = '  +5.47%  '
$ws.Range("D12").Value = '0.0805'
$ws.Range("E12").Value = '  +0.29%  '
$ws.Range("D13").Value = '0.113'
$ws.Range("E13").Value = '  -2.91%  '
$ws.Range("D14").Value = '6.73'
$ws.Range("E14").Value = '  +0.31%  '
$ws.Range("D15").Value = '2.638.45'
$ws.Range("E15").Value = '  +0.36%  '
$ws.Range("D16").Value = '14.42'
$ws.Range("E16").Value = '  +0.85%  '
$ws.Range("D17").Value = '2.294.33'
$ws.Range("E17").Value = '  +0.17%  '
$ws.Range("D18").Value = '0.792'
$ws.Range("E18").Value = '  +4.10%  '
$ws.Range("D19").Value = '42.133.40'
$ws.Range("E19").Value = '  +0.85%  '
$ws.Range("D20").Value = '12.89'
$ws.Range("E20").Value = '  +4.90%  '
$ws.Range("D21").Value = '0.0₃0921'
$ws.Range("E21").Value = '  +1.35%  '
$ws.Range("D22").Value = '6.03'
$ws.Range("E22").Value = '  +0.96%  '
$ws.Range("D23").Value = '68.21'
$ws.Range("E23").Value = '  +1.14%  '
$ws.Range("D24").Value = '244.56'
$ws.Range("E24").Value = '  +0.29%  '
$ws.Range("D25").Value = '2.61'
$ws.Range("E25").Value = '  +0.83%  '
$ws.Range("D26").Value = '1.96'
$ws.Range("E26").Value = '  +1.06%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").Value = '24.17'
$ws.Range("E28").Value = '  -1.12%  '
$ws.Range("D29").Value = '9.73'
$ws.Range("E29").Value = '  +0.81%  '
$ws.Range("D30").Value = '35.69'
$ws.Range("E30").Value = '  +4.02%  '
$ws.Range("D31").Value = '2.10'
$ws.Range("E31").Value = '  +1.18%  '
$ws.Range("D32").Value = '160.80'
$ws.Range("E32").Value = '  +1.47%  '
$ws.Range("D33").Value = '5.37'
$ws.Range("E33").Value = '  +3.22%  '
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  +0.07%  '
$ws.Range("D35").Value = '0.0753'
$ws.Range("E35").Value = '  +0.38%  '
$ws.Range("D36").Value = '3.11'
$ws.Range("E36").Value = '  +0.91%  '
$ws.Range("D37").Value = '17.33'
$ws.Range("E37").Value = '  +2.75%  '
$ws.Range("E38").Value = '  +2.22%  '
$ws.Range("D39").Value = '2.38'
$ws.Range("E39").Value = '  -0.47%  '
$ws.Range("E40").Value = '  +0.14%  '
$ws.Range("D41").Value = '1.82'
$ws.Range("E41").Value = '  -0.69%  '
$ws.Range("D42").Value = '4.16'
$ws.Range("E42").Value = '  +5.35%  '
$ws.Range("D43").Value = '19.94'
$ws.Range("E43").Value = '  +0.72%  '
$ws.Range("D44").Value = '2.017.99'
$ws.Range("E44").Value = '  -2.73%  '
$ws.Range("D45").Value = '2.27'
$ws.Range("E45").Value = '  +11.10%  '
$ws.Range("D46").Value = '0.0284'
$ws.Range("E46").Value = '  +1.24%  '
$ws.Range("D47").Value = '10.26'
$ws.Range("E47").Value = '  -1.64%  '
$ws.Range("D48").Value = '2.96'
$ws.Range("E48").Value = '  +0.87%  '
$ws.Range("D49").Value = '53.73'
$ws.Range("E49").Value = '  +3.12%  '
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").Value = '1.53'
$ws.Range("E50").Value = '  -0.69%  '
$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D51").Value = '72.60'
$ws.Range("E51").Value = '  -1.02%  '
